# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.413.15"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.721.01"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5316"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06716"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2669"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.02"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.517"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "1.957.63"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "1.721.33"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5868"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "0.0₅8232"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.26"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "27.457.05"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "224.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.678"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.52"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.064"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.34"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.700"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1212"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.274"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05391"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.296"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.494"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.433"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.633"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.871"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9570"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.391"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5907"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").Value = "1.157.41"
$ws.Range("E39").Value = "  +10.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01656"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8436"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.92"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "1.864.30"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -6.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.04"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4585"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.195"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05210"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.58%  "
